$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 658, shifting existing rows 658..693 down to 659..694.
$ws.Rows.Item(658).Insert()

# Fill the new row 658 with its data. Columns A, B, C, E, F, G, I, R repeat the
# same constant values used throughout this block of rows.
$ws.Cells.Item(658, 1).Value = 5
$ws.Cells.Item(658, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(658, 3).Value = "Maule"
$ws.Cells.Item(658, 4).Value = 45267
$ws.Cells.Item(658, 5).Value = 7
$ws.Cells.Item(658, 6).Value = 100112032
$ws.Cells.Item(658, 7).Value = "Zapallo italiano"
$ws.Cells.Item(658, 8).Value = "Sin especificar"
$ws.Cells.Item(658, 9).Value = "Primera"
$ws.Cells.Item(658, 10).Value = 500
$ws.Cells.Item(658, 11).Value = 7000
$ws.Cells.Item(658, 12).Value = 7000
$ws.Cells.Item(658, 13).Value = 7000
$ws.Cells.Item(658, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(658, 15).Value = "Región del Maule"
$ws.Cells.Item(658, 16).Value = 140
$ws.Cells.Item(658, 17).Value = 50
$ws.Cells.Item(658, 18).Value = "Hortaliza"

# Apply the same date-number-format style used by the other rows' date column
# (D) to the newly inserted row's D cell.
$ws.Cells.Item(658, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
